$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "3+82="
$t.Cell(1,2).Range.Text = "83-74="
$t.Cell(1,3).Range.Text = "50-12="
$t.Cell(1,4).Range.Text = "83-50="
$t.Cell(1,5).Range.Text = "27+70="
$t.Cell(2,1).Range.Text = "38+61="
$t.Cell(2,2).Range.Text = "63+5="
$t.Cell(2,3).Range.Text = "72+6="
$t.Cell(2,4).Range.Text = "80-33="
$t.Cell(2,5).Range.Text = "93-11="
$t.Cell(3,1).Range.Text = "17+15="
$t.Cell(3,2).Range.Text = "44-40="
$t.Cell(3,3).Range.Text = "30+32="
$t.Cell(3,4).Range.Text = "23+20="
$t.Cell(3,5).Range.Text = "84-21="
$t.Cell(4,1).Range.Text = "36+45="
$t.Cell(4,2).Range.Text = "9+56="
$t.Cell(4,3).Range.Text = "0+89="
$t.Cell(4,4).Range.Text = "57+11="
$t.Cell(4,5).Range.Text = "94-9="
$t.Cell(5,1).Range.Text = "48-14="
$t.Cell(5,2).Range.Text = "58-46="
$t.Cell(5,3).Range.Text = "2+24="
$t.Cell(5,4).Range.Text = "17+6="
$t.Cell(5,5).Range.Text = "43+41="
$t.Cell(6,1).Range.Text = "48+38="
$t.Cell(6,2).Range.Text = "25+36="
$t.Cell(6,3).Range.Text = "96-32="
$t.Cell(6,4).Range.Text = "87-42="
$t.Cell(6,5).Range.Text = "26+24="
$t.Cell(7,1).Range.Text = "34-26="
$t.Cell(7,2).Range.Text = "49-35="
$t.Cell(7,3).Range.Text = "90-12="
$t.Cell(7,4).Range.Text = "90-84="
$t.Cell(7,5).Range.Text = "84+12="
$t.Cell(8,1).Range.Text = "66-64="
$t.Cell(8,2).Range.Text = "85-72="
$t.Cell(8,3).Range.Text = "3+54="
$t.Cell(8,4).Range.Text = "51-12="
$t.Cell(8,5).Range.Text = "42+49="
$t.Cell(9,1).Range.Text = "16+40="
$t.Cell(9,2).Range.Text = "2+38="
$t.Cell(9,3).Range.Text = "32+40="
$t.Cell(9,4).Range.Text = "50+4="
$t.Cell(9,5).Range.Text = "79+17="
$t.Cell(10,1).Range.Text = "93-84="
$t.Cell(10,2).Range.Text = "54-48="
$t.Cell(10,3).Range.Text = "11+80="
$t.Cell(10,4).Range.Text = "39-9="
$t.Cell(10,5).Range.Text = "39+10="
$t.Cell(11,1).Range.Text = "3+75="
$t.Cell(11,2).Range.Text = "52+15="
$t.Cell(11,3).Range.Text = "14+13="
$t.Cell(11,4).Range.Text = "12+70="
$t.Cell(11,5).Range.Text = "69-59="
$t.Cell(12,1).Range.Text = "74-4="
$t.Cell(12,2).Range.Text = "71-43="
$t.Cell(12,3).Range.Text = "69-39="
$t.Cell(12,4).Range.Text = "46-7="
$t.Cell(12,5).Range.Text = "3+33="
$t.Cell(13,1).Range.Text = "68-11="
$t.Cell(13,2).Range.Text = "59-32="
$t.Cell(13,3).Range.Text = "78-46="
$t.Cell(13,4).Range.Text = "4+66="
$t.Cell(13,5).Range.Text = "17+62="
$t.Cell(14,1).Range.Text = "26+25="
$t.Cell(14,2).Range.Text = "50-17="
$t.Cell(14,3).Range.Text = "61+30="
$t.Cell(14,4).Range.Text = "68+5="
$t.Cell(14,5).Range.Text = "83-5="
$t.Cell(15,1).Range.Text = "87-62="
$t.Cell(15,2).Range.Text = "80-50="
$t.Cell(15,3).Range.Text = "30+25="
$t.Cell(15,4).Range.Text = "86-26="
$t.Cell(15,5).Range.Text = "88+0="
$t.Cell(16,1).Range.Text = "68-40="
$t.Cell(16,2).Range.Text = "99-88="
$t.Cell(16,3).Range.Text = "25+29="
$t.Cell(16,4).Range.Text = "86+7="
$t.Cell(16,5).Range.Text = "6+65="
$t.Cell(17,1).Range.Text = "63+29="
$t.Cell(17,2).Range.Text = "93-65="
$t.Cell(17,3).Range.Text = "5+70="
$t.Cell(17,4).Range.Text = "31-30="
$t.Cell(17,5).Range.Text = "91+0="
$t.Cell(18,1).Range.Text = "79-57="
$t.Cell(18,2).Range.Text = "21+50="
$t.Cell(18,3).Range.Text = "51-23="
$t.Cell(18,4).Range.Text = "99-77="
$t.Cell(18,5).Range.Text = "14+57="
$t.Cell(19,1).Range.Text = "18-18="
$t.Cell(19,2).Range.Text = "33-6="
$t.Cell(19,3).Range.Text = "36-25="
$t.Cell(19,4).Range.Text = "25+60="
$t.Cell(19,5).Range.Text = "33+4="
$t.Cell(20,1).Range.Text = "58+13="
$t.Cell(20,2).Range.Text = "10+55="
$t.Cell(20,3).Range.Text = "79-20="
$t.Cell(20,4).Range.Text = "18+16="
$t.Cell(20,5).Range.Text = "6+17="
